# Apply hybrid bold + color (2C3E50) highlighting to quantitative metrics
# (percentages, dollar amounts, large numbers) across achievement / work
# experience bullet paragraphs, per the commit:
#   "Implement quantitative metrics highlighting across all resume formats"

function Bold-Metrics {
    param($doc, $paraIndex, $substrings)

    $p = $doc.Paragraphs.Item($paraIndex)
    $pStart = $p.Range.Start
    $pEnd = $p.Range.End
    $cursor = $pStart

    foreach ($t in $substrings) {
        $searchRng = $doc.Range($cursor, $pEnd)
        $f = $searchRng.Find.Execute($t, $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
        if (-not $f) {
            Write-Output ("WARNING: substring not found -> para " + $paraIndex + " [" + $t + "]")
        } else {
            $searchRng.Font.Bold = 1
            $searchRng.Font.Color = 5258796
            $cursor = $searchRng.End
        }
    }
}

$d = $word.ActiveDocument

# Partner - Siege Analytics bullets
Bold-Metrics $d 10 @("23%", "64%")
Bold-Metrics $d 12 @([char]0x00B1 + "4.2%", [char]0x00B1 + "2.1%", "71%", "87%")
Bold-Metrics $d 13 @("73.5%", "`$4.7M")
Bold-Metrics $d 14 @("`$2")

# Data Products Manager - Helm/Murmuration bullet
Bold-Metrics $d 20 @("57%")

# KEY ACHIEVEMENTS AND IMPACT bullets
Bold-Metrics $d 85 @("178%")
Bold-Metrics $d 86 @("73.5%")
Bold-Metrics $d 87 @("`$4.7M")
Bold-Metrics $d 88 @("12,847")
Bold-Metrics $d 90 @([char]0x00B1 + "4.2%", [char]0x00B1 + "2.1%")
Bold-Metrics $d 91 @("71%", "87%")

Write-Output "Metrics highlighting applied."
